$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# George paid his June subscription (row 10 is George's row, F column = June,18)
$ws.Range("F10").Value = 300

# New booking history entry: slot booked for 26/06/2018, 4-5, balance updated accordingly
$ws.Range("B52").Value = "26/06/2018, 4-5"

$ws.Range("C52").Value2 = 43277
$ws.Range("C41").Copy()
$ws.Range("C52").PasteSpecial(-4122)  # xlPasteFormats, reuse existing date-format style

$ws.Range("D52").Value = 1265

$excel.CutCopyMode = 0

# Scroll/selection state left by the author after the edit
$ws.Application.ActiveWindow.ScrollRow = 31
$ws.Range("D52").Select()
